$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '25.930.86'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.35%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.615.42'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -1.14%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '211.63'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.17%  '
$ws.Range('E6').Value = '  +0.02%  '
$ws.Range('E7').Value = '  -3.61%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.0621'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.45%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.246'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.89%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '18.14'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.39%  '
$ws.Range('E11').Value = '  -0.40%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.841.36'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.06%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.629.64'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.41%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.10'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.61%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.516'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -2.66%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '25.941.14'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.29%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '61.46'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.78%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.0₃0730'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -2.25%  '
$ws.Range('E19').Value = '  +0.03%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '190.76'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.02%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.21'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.21%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.42'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.73%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.00'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.41%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.130'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.42%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '143.33'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.03%  '
$ws.Range('E26').Value = '  +0.04%  '
$ws.Range('E27').Value = '  -3.35%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.59'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -2.56%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.12'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.80%  '
$ws.Range('E30').Value = '  -1.12%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0472'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -2.56%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.10'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.75%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.07'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.88%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.41'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.03%  '
$ws.Range('E35').Value = '  -1.54%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.123.01'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.74%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.816'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -6.77%  '
$ws.Range('E38').Value = '  -1.79%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.513'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.32%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0152'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.83%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '97.36'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.35%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.752.60'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.01%  '
$ws.Range('E43').Value = '  -4.05%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.07'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -4.36%  '
$ws.Range('E45').Value = '  -0.28%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '53.66'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.99%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.47'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.19%  '
$ws.Range('E48').Value = '  -2.36%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.411'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.73%  '
$ws.Range('E50').Value = '  -0.17%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.38'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.23%  '
